# Applies the Halicarnassus_Profits price/profit refresh produced by the
# scheduled Sheets runner: updates currentAveragePrice* / Leve* columns
# (H, I, J, K, L) and the derived LeveProfit columns (M, N) for the rows
# whose source data changed, across all eight crafting-job worksheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 440
$ws.Range("I21").Value = 440
$ws.Range("K21").Value = 440
$ws.Range("M21").Value = 28
$ws.Range("H23").Value = 440
$ws.Range("I23").Value = 440
$ws.Range("K23").Value = 440
$ws.Range("M23").Value = -206
$ws.Range("H32").Value = 624.75
$ws.Range("J32").Value = 749.5
$ws.Range("L32").Value = 749.5
$ws.Range("N32").Value = -1401.5
$ws.Range("H34").Value = 21250
$ws.Range("I34").Value = 21250
$ws.Range("K34").Value = 21250
$ws.Range("M34").Value = -21047
$ws.Range("H36").Value = 21250
$ws.Range("I36").Value = 21250
$ws.Range("K36").Value = 21250
$ws.Range("M36").Value = -20535
$ws.Range("H43").Value = 10002419
$ws.Range("I43").Value = 12502024
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 12502024
$ws.Range("L43").Value = 4000
$ws.Range("M43").Value = -12501955
$ws.Range("N43").Value = -4138
$ws.Range("H49").Value = 508.5
$ws.Range("I49").Value = 17
$ws.Range("K49").Value = 51
$ws.Range("M49").Value = 85
$ws.Range("H64").Value = 7816.6665
$ws.Range("I64").Value = 5950
$ws.Range("K64").Value = 5950
$ws.Range("M64").Value = -5702
$ws.Range("H67").Value = 7816.6665
$ws.Range("I67").Value = 5950
$ws.Range("K67").Value = 5950
$ws.Range("M67").Value = -5092
$ws.Range("H98").Value = 1058.4
$ws.Range("I98").Value = 1058.4
$ws.Range("K98").Value = 1058.4
$ws.Range("M98").Value = 439.5999999999999
$ws.Range("H122").Value = 1058.4
$ws.Range("I122").Value = 1058.4
$ws.Range("K122").Value = 3175.2
$ws.Range("M122").Value = -725.2000000000003
$ws.Range("H137").Value = 3120.647
$ws.Range("J137").Value = 3405.4285
$ws.Range("L137").Value = 10216.2855
$ws.Range("N137").Value = -15316.2855

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 30000
$ws.Range("I31").Value = 30000
$ws.Range("K31").Value = 30000
$ws.Range("M31").Value = -29706
$ws.Range("H74").Value = 1779.6
$ws.Range("I74").Value = 1779.6
$ws.Range("K74").Value = 1779.6
$ws.Range("M74").Value = -905.5999999999999
$ws.Range("H77").Value = 1779.6
$ws.Range("I77").Value = 1779.6
$ws.Range("K77").Value = 8898
$ws.Range("M77").Value = -4530
$ws.Range("H122").Value = 3485.3333
$ws.Range("I122").Value = 3512.7144
$ws.Range("J122").Value = 3102
$ws.Range("K122").Value = 10538.1432
$ws.Range("L122").Value = 9306
$ws.Range("M122").Value = -8088.143199999999
$ws.Range("N122").Value = -14206

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 713.25
$ws.Range("I22").Value = 557
$ws.Range("K22").Value = 557
$ws.Range("M22").Value = -384
$ws.Range("H80").Value = 608.4286
$ws.Range("I80").Value = 359.6
$ws.Range("K80").Value = 359.6
$ws.Range("M80").Value = 638.4
$ws.Range("H83").Value = 608.4286
$ws.Range("I83").Value = 359.6
$ws.Range("K83").Value = 1798
$ws.Range("M83").Value = 3194
$ws.Range("H105").Value = 1693.75
$ws.Range("I105").Value = 1650
$ws.Range("J105").Value = 1766.6666
$ws.Range("K105").Value = 1650
$ws.Range("L105").Value = 1766.6666
$ws.Range("M105").Value = 97
$ws.Range("N105").Value = -5260.6666

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1821.3572
$ws.Range("I58").Value = 964.0909
$ws.Range("J58").Value = 4964.6665
$ws.Range("K58").Value = 964.0909
$ws.Range("L58").Value = 4964.6665
$ws.Range("M58").Value = -761.0909
$ws.Range("N58").Value = -5370.6665
$ws.Range("H136").Value = 1821.3572
$ws.Range("I136").Value = 964.0909
$ws.Range("J136").Value = 4964.6665
$ws.Range("K136").Value = 2892.2727
$ws.Range("L136").Value = 14893.9995
$ws.Range("M136").Value = -342.2727
$ws.Range("N136").Value = -19993.9995

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2119.6667
$ws.Range("I109").Value = 2119.6667
$ws.Range("K109").Value = 6359.000100000001
$ws.Range("M109").Value = -5319.000100000001
$ws.Range("H122").Value = 900
$ws.Range("I122").Value = 900
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8100
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -5650
$ws.Range("H132").Value = 2135.2856
$ws.Range("I132").Value = 2082.6667
$ws.Range("J132").Value = 2174.75
$ws.Range("K132").Value = 18744.0003
$ws.Range("L132").Value = 19572.75
$ws.Range("M132").Value = -16214.0003
$ws.Range("N132").Value = -24632.75
$ws.Range("H136").Value = 7559.1665
$ws.Range("I136").Value = 6980
$ws.Range("K136").Value = 20940
$ws.Range("M136").Value = -15840

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2513
$ws.Range("I102").Value = 2071.4443
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 2071.4443
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = -449.4443000000001
$ws.Range("N102").Value = -7744

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4005.75
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -505
$ws.Range("H27").Value = 4005.75
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 800
$ws.Range("M27").Value = -693
$ws.Range("H40").Value = 1749.25
$ws.Range("I40").Value = 1749.25
$ws.Range("K40").Value = 1749.25
$ws.Range("M40").Value = -1613.25
$ws.Range("H46").Value = 5436.385
$ws.Range("J46").Value = 6374.8335
$ws.Range("L46").Value = 6374.8335
$ws.Range("N46").Value = -6750.8335
$ws.Range("H55").Value = 4751.231
$ws.Range("I55").Value = 4922.6665
$ws.Range("J55").Value = 4699.8
$ws.Range("K55").Value = 4922.6665
$ws.Range("L55").Value = 4699.8
$ws.Range("M55").Value = -4749.6665
$ws.Range("N55").Value = -5045.8
$ws.Range("H93").Value = 5181.273
$ws.Range("I93").Value = 2331.3333
$ws.Range("K93").Value = 2331.3333
$ws.Range("M93").Value = -1083.3333

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1963.7273
$ws.Range("I122").Value = 1860.1
$ws.Range("K122").Value = 5580.299999999999
$ws.Range("M122").Value = -3130.299999999999

